$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.352819
$ws.Range("H2").Value = 10.058457
$ws.Range("I2").Value = 0.02224149976981271
$ws.Range("J2").Value = 0.02224149976981271
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.235341333333333
$ws.Range("N2").Value = 9.706023999999999
$ws.Range("O2").Value = 0.2153734454473681
$ws.Range("P2").Value = 0.2153734454473681
$ws.Range("Q2").Value = 10.84751389388533
$ws.Range("R2").Value = 97.627625044968
$ws.Range("S2").Value = 0.004790228437341407
$ws.Range("T2").Value = 0.004790228437341408
$ws.Range("G3").Value = 3.352819
$ws.Range("H3").Value = 10.058457
$ws.Range("I3").Value = 0.02224149976981271
$ws.Range("J3").Value = 0.02224149976981271
$ws.Range("O3").Value = 0.4841904166376352
$ws.Range("P3").Value = 0.4841904166376352
$ws.Range("Q3").Value = 24.38676811272167
$ws.Range("R3").Value = 219.480913014495
$ws.Range("S3").Value = 0.01076912104019148
$ws.Range("T3").Value = 0.01076912104019148
$ws.Range("G4").Value = 3.352819
$ws.Range("H4").Value = 10.058457
$ws.Range("I4").Value = 0.02224149976981271
$ws.Range("J4").Value = 0.02224149976981271
$ws.Range("O4").Value = 0.3004361379149967
$ws.Range("P4").Value = 0.3004361379149967
$ws.Range("Q4").Value = 15.13178736351967
$ws.Range("R4").Value = 136.186086271677
$ws.Range("S4").Value = 0.006682150292279819
$ws.Range("T4").Value = 0.006682150292279819
$ws.Range("I5").Value = 0.8292884613633072
$ws.Range("J5").Value = 0.8292884613633072
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.235341333333333
$ws.Range("N5").Value = 9.706023999999999
$ws.Range("O5").Value = 0.2153734454473681
$ws.Range("P5").Value = 0.2153734454473681
$ws.Range("Q5").Value = 404.4564530170178
$ws.Range("R5").Value = 3640.10807715316
$ws.Range("S5").Value = 0.1786067131935621
$ws.Range("T5").Value = 0.1786067131935621
$ws.Range("I6").Value = 0.8292884613633072
$ws.Range("J6").Value = 0.8292884613633072
$ws.Range("O6").Value = 0.4841904166376352
$ws.Range("P6").Value = 0.4841904166376352
$ws.Range("Q6").Value = 909.2761556156973
$ws.Range("R6").Value = 8183.485400541275
$ws.Range("S6").Value = 0.4015335256202832
$ws.Range("T6").Value = 0.4015335256202832
$ws.Range("I7").Value = 0.8292884613633072
$ws.Range("J7").Value = 0.8292884613633072
$ws.Range("O7").Value = 0.3004361379149967
$ws.Range("P7").Value = 0.3004361379149967
$ws.Range("Q7").Value = 564.1983135238739
$ws.Range("S7").Value = 0.249148222549462
$ws.Range("T7").Value = 0.249148222549462
$ws.Range("I8").Value = 0.1484700388668802
$ws.Range("J8").Value = 0.1484700388668802
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3.235341333333333
$ws.Range("N8").Value = 9.706023999999999
$ws.Range("O8").Value = 0.2153734454473681
$ws.Range("P8").Value = 0.2153734454473681
$ws.Range("Q8").Value = 72.41107057088264
$ws.Range("R8").Value = 651.6996351379439
$ws.Range("S8").Value = 0.03197650381646464
$ws.Range("T8").Value = 0.03197650381646464
$ws.Range("I9").Value = 0.1484700388668802
$ws.Range("J9").Value = 0.1484700388668802
$ws.Range("O9").Value = 0.4841904166376352
$ws.Range("P9").Value = 0.4841904166376352
$ws.Range("S9").Value = 0.07188776997716061
$ws.Range("T9").Value = 0.07188776997716061
$ws.Range("I10").Value = 0.1484700388668802
$ws.Range("J10").Value = 0.1484700388668802
$ws.Range("O10").Value = 0.3004361379149967
$ws.Range("P10").Value = 0.3004361379149967
$ws.Range("R10").Value = 909.0912812150908
$ws.Range("S10").Value = 0.04460576507325494
$ws.Range("T10").Value = 0.04460576507325494
